# Apply updated cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '49.423.11'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.628.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '111.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '325.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.524'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.549'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.77%  '
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0811'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.040.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.644.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.853'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '49.389.00'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0947'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '268.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.20%  '
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.137'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.39'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.56'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0806'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.09'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '129.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.36'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.76'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.31%  '
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0335'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.061.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('E47').Value = '  +7.78%  '
$ws.Range('E48').Value = '  -6.86%  '
$ws.Range('E49').Value = '  -3.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '58.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.85%  '
